$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 'switzerland; usa; italy; philipines; uk; luxembourg; kenya; sweden; canada; netherlands; south africa; germany'
$ws.Range("C4").Value = 'usa; canada; uk; brazil; australia; italy; china; slovakia; spain; spain; sweden; new zealand'
$ws.Range("D2").Value = 'usa; uk; india; australia; france; bangladesh; brazil; germany; philipines; new zealand; russia; canada; israel; italy; japan; malysia; netherlands; singapor; south africa; ukraine; uae; armenia; belgium; croatia; estonia; ireland; jamaica; jordan; nepal; nigeria; portugal; sweden; switzerland; turkey'
$ws.Range("D4").Value = 'usa; uk; australia; india; germany; spain; canada; netherlands; new zealand; italy; singapor; argentina'
$ws.Range("D9").Value = 'usa; uk; india; australia; germany; spain; japan; russia; new zealand; singapor; canada; france'
$ws.Range("D13").Value = 'usa; uk; india; australia; germany; netherlands; singapor; canada; spain; france; italy '
$ws.Range("B16").Value = 'uk; usa; singapor; india; brunei; france; australia; finland; germany; hungary; iran'
$ws.Range("C2").Value = 'usa; uk; australia; italy; costa rica; germany; greece; ireland; mexico; new zealand; sweden'
$ws.Range("C3").Value = 'usa; uk; canada; brazil; italy; india; australia; spain; norway; serbia; thailand'
$ws.Range("D3").Value = 'usa; uk; australia; india; germany; spain; canada; russia; france; italy; iran'
$ws.Range("D15").Value = 'usa; uk; india; australia; germany; canada; spain; south africa; argentina; france; netherlands'
$ws.Range("E7").Value = 'switzerland; uk; luxembourg; canada; italy; netherlands; usa; belgium; chile; south africa; sweden; uae'
$ws.Range("C6").Value = 'usa; uk; canada; brazil; australia; south africa; greece; india; italy; japan; new zealand'
$ws.Range("B10").Value = 'usa; brazil; uk; canada; france; germany; spain; iran; australia; japan; switzerland'
$ws.Range("B16").Select()
